# Update the QuantitativeMetrics evaluation sheet with the latest test run
# results (test code generation module - update evaluations).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")
$ws.Activate()

# Assertion validity: the test stopped passing; note a portal issue instead
# of the previous baseline comment.
$ws.Range("B7").Value = "no"
$ws.Range("C7").Value = "Portal issue"

# Updated Code BLEU score (dataflow_match_score dropped, the rest unchanged).
$ws.Range("B12").Value = 0.3070794321608488
$ws.Range("C12").Value = "{'codebleu': 0.30707943216084876, 'ngram_match_score': 0.0657276363698262, 'weighted_ngram_match_score': 0.0812297823295309, 'syntax_match_score': 0.6683168316831684, 'dataflow_match_score': 0.41304347826086957}"

# Leave the selection where the author last left it when saving.
$ws.Range("B8").Select()
